# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracking sheet.
# The sheet keeps a special "date only" number format on the A-cell of the
# last row, while every other date cell uses the full date-time format.
# So before adding the new last row we need to reset the previous last
# row (A55) back to the standard format, then give the new last row
# (A56) the "last row" date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the previous last row's date cell to the regular (non-final) format.
$ws.Range("A55").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new day's data.
$ws.Range("A56").Value = 45796
$ws.Range("B56").Value = 232
$ws.Range("C56").Value = 242
$ws.Range("D56").Value = 236

# Give the new last row's date cell the "final row" date-only format.
$ws.Range("A56").NumberFormat = "YYYY-MM-DD"
